# Replace a unique run of plain text with new text while avoiding Word's
# automatic coalescing of the edited run with its (untouched) neighboring
# runs. We do this by temporarily splitting the paragraph right after the
# matched text (so the matched run becomes the last run of a throw-away
# paragraph), performing a Delete()+InsertAfter() there -- which does not
# pull neighboring runs into the edit -- and then removing the temporary
# paragraph mark to rejoin the text back into a single paragraph.
function Replace-Text($d, $oldText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $e = $rng.End
    $splitPoint = $d.Range($e, $e)
    $splitPoint.InsertParagraphAfter()

    $rng2 = $d.Content
    $rng2.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng2.Delete()
    $rng2.InsertAfter($newText)

    $e2 = $rng2.End
    $markRng = $d.Range($e2, $e2 + 1)
    $markRng.Delete()
}

$d = $word.ActiveDocument

# "Fit the data: Using Scikit Learn Library" block (appears twice in the
# document with identical before/after numbers).
Replace-Text $d "0.9863 and coefficients" "0.9776 and coefficients"
Replace-Text $d "0.9863 and coefficients" "0.9776 and coefficients"

Replace-Text $d "3.0371, and" "3.0198, and"
Replace-Text $d "3.0371, and" "3.0198, and"

Replace-Text $d "1.9549" "2.0666"
Replace-Text $d "1.9549" "2.0666"

# "Fit the data: Using Gradient Descent" block.
Replace-Text $d "0.9857 and coefficients" "0.9771 and coefficients"
Replace-Text $d "3.0359, and" "3.0184, and"
Replace-Text $d "1.9543" "2.066"

# Custom-library / array(...) block.
Replace-Text $d "array([0.99955242]) and coefficients" "array([0.98750421]) and coefficients"
Replace-Text $d "array([3.0286027]), and" "array([3.01288453]), and"
Replace-Text $d "array([1.93519378])" "array([2.05965667])"
